$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, shifting existing rows 4..73 down to 5..74
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112052
$ws.Cells.Item(4, 7).Value = "Albahaca"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 900
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 13).Value = 950
$ws.Cells.Item(4, 14).Value = "$/paquete"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 950
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Match the date formatting/style used by the other rows in column D
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
